$wb = $excel.ActiveWorkbook

# --- Cell value edits: browser column entries changed to "chrome" ---
# Sheet1 (verifyAvailablePets): A2 was "edge" -> "chrome"
$wb.Worksheets.Item(1).Range("A2").Value = "chrome"

# Sheet4 (VerifySignInButtonPresence): A2 was "firefox" -> "chrome"
$wb.Worksheets.Item(4).Range("A2").Value = "chrome"

# Sheet5 (Verifylogin): A3 and A5 were "firefox" -> "chrome"
$wb.Worksheets.Item(5).Range("A3").Value = "chrome"
$wb.Worksheets.Item(5).Range("A5").Value = "chrome"

# --- Update selections / active cells on each sheet ---
[void]$wb.Worksheets.Item(4).Activate()
[void]$wb.Worksheets.Item(4).Range("A11").Select()

[void]$wb.Worksheets.Item(5).Activate()
[void]$wb.Worksheets.Item(5).Range("A6").Select()

[void]$wb.Worksheets.Item(2).Activate()
[void]$wb.Worksheets.Item(2).Range("J20").Select()

[void]$wb.Worksheets.Item(3).Activate()
[void]$wb.Worksheets.Item(3).Range("I25").Select()

[void]$wb.Worksheets.Item(6).Activate()
[void]$wb.Worksheets.Item(6).Range("G8").Select()

# Sheet1 ends up as the active/selected sheet
[void]$wb.Worksheets.Item(1).Activate()
[void]$wb.Worksheets.Item(1).Range("L22").Select()
